# Updates crypto price/volume data (and the Filecoin/Stacks row-content swap at rows 41-42)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.270.43"
$ws.Range("E2").Value = "'  -1.43%  "
$ws.Range("D3").Value = "'3.076.67"
$ws.Range("E3").Value = "'  -1.31%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'522.01"
$ws.Range("E5").Value = "'  -1.08%  "
$ws.Range("D6").Value = "'135.70"
$ws.Range("E6").Value = "'  -4.67%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  -0.08%  "
$ws.Range("D8").Value = "'3.077.98"
$ws.Range("E8").Value = "'  -1.23%  "
$ws.Range("E9").Value = "'  +4.48%  "
$ws.Range("D10").Value = "'7.29"
$ws.Range("E10").Value = "'  +1.81%  "
$ws.Range("E11").Value = "'  -2.05%  "
$ws.Range("D12").Value = "'0.400"
$ws.Range("E12").Value = "'  +1.78%  "
$ws.Range("E13").Value = "'  +1.27%  "
$ws.Range("D14").Value = "'3.607.72"
$ws.Range("E14").Value = "'  -1.27%  "
$ws.Range("D15").Value = "'25.22"
$ws.Range("E15").Value = "'  -1.52%  "
$ws.Range("E16").Value = "'  -2.18%  "
$ws.Range("D17").Value = "'57.353.95"
$ws.Range("E17").Value = "'  -1.33%  "
$ws.Range("D18").Value = "'3.070.04"
$ws.Range("E18").Value = "'  -1.43%  "
$ws.Range("E19").Value = "'  -4.15%  "
$ws.Range("D21").Value = "'7.82"
$ws.Range("E21").Value = "'  -1.89%  "
$ws.Range("D22").Value = "'350.89"
$ws.Range("E22").Value = "'  +2.51%  "
$ws.Range("E23").Value = "'  +0.10%  "
$ws.Range("D24").Value = "'69.02"
$ws.Range("E24").Value = "'  +2.13%  "
$ws.Range("E25").Value = "'  -3.05%  "
$ws.Range("E26").Value = "'  -2.66%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "'  -0.37%  "
$ws.Range("D28").Value = "'0.0₃0867"
$ws.Range("E28").Value = "'  -6.07%  "
$ws.Range("E29").Value = "'  +0.04%  "
$ws.Range("D30").Value = "'7.21"
$ws.Range("E30").Value = "'  -1.60%  "
$ws.Range("E31").Value = "'  -0.90%  "
$ws.Range("E32").Value = "'  -8.59%  "
$ws.Range("D33").Value = "'20.92"
$ws.Range("E33").Value = "'  -0.57%  "
$ws.Range("D34").Value = "'4.85"
$ws.Range("E34").Value = "'  +3.19%  "
$ws.Range("D35").Value = "'159.16"
$ws.Range("E35").Value = "'  +0.46%  "
$ws.Range("E36").Value = "'  -4.81%  "
$ws.Range("D37").Value = "'6.01"
$ws.Range("E37").Value = "'  -3.23%  "
$ws.Range("D38").Value = "'25.43"
$ws.Range("E38").Value = "'  -3.59%  "
$ws.Range("E39").Value = "'  -1.58%  "
$ws.Range("E40").Value = "'  -1.51%  "
$ws.Range("B41").Value = "'Filecoin"
$ws.Range("C41").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'4.06"
$ws.Range("E41").Value = "'  +1.90%  "
$ws.Range("B42").Value = "'Stacks"
$ws.Range("C42").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.57"
$ws.Range("E42").Value = "'  -5.51%  "
$ws.Range("D43").Value = "'0.693"
$ws.Range("E43").Value = "'  +0.06%  "
$ws.Range("D44").Value = "'2.408.56"
$ws.Range("E44").Value = "'  +5.86%  "
$ws.Range("E45").Value = "'  -0.07%  "
$ws.Range("E46").Value = "'  +0.03%  "
$ws.Range("D47").Value = "'3.117.74"
$ws.Range("E47").Value = "'  -1.27%  "
$ws.Range("E48").Value = "'  -0.08%  "
$ws.Range("D49").Value = "'0.946"
$ws.Range("E49").Value = "'  -5.35%  "
$ws.Range("E50").Value = "'  -2.29%  "
$ws.Range("E51").Value = "'  -5.31%  "
